# Add GFLOP/sec column (E) to the three "strong overlap cpu" sheets,
# then add a new GFLOP/sec strong-scaling chart on the summary sheet,
# matching the commit "Added a new gflops figure for strong scaling".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. n50_strong_overlap_cpu : E1 header + E2:E12 GFLOP/sec formulas
# ---------------------------------------------------------------
$ws50 = $wb.Worksheets.Item("n50_strong_overlap_cpu")
$ws50.Range("E1").Value = "GFLOPs"
$ws50.Range("E2").Formula = "=((4096000*50*2)/(B2*0.001))*0.000000001"
$ws50.Range("E3").Formula = "=((4096000*50*2)/(B3*0.001))*0.000000001"
$ws50.Range("E4:E12").Formula = "=((4096000*50*2)/(B4*0.001))*0.000000001"

# ---------------------------------------------------------------
# 2. n17_strong_overlap_cpu : E1 header + E2:E12 GFLOP/sec formulas
# ---------------------------------------------------------------
$ws17 = $wb.Worksheets.Item("n17_strong_overlap_cpu")
$ws17.Range("E1").Value = "GFLOPs"
$ws17.Range("E2:E11").Formula = "=((4096000*17*2)/(B2*0.001))*0.000000001"
$ws17.Range("E12").Formula = "=((4096000*17*2)/(B12*0.001))*0.000000001"

# ---------------------------------------------------------------
# 3. n31_strong_overlap_cpu : E1 header + E2:E12 GFLOP/sec formulas
# ---------------------------------------------------------------
$ws31 = $wb.Worksheets.Item("n31_strong_overlap_cpu")
$ws31.Range("E1").Value = "GFLOPs"
$ws31.Range("E2").Formula = "=((4096000*31*2)/(B2*0.001))*0.000000001"
$ws31.Range("E3:E12").Formula = "=((4096000*31*2)/(B3*0.001))*0.000000001"

# ---------------------------------------------------------------
# 4. New chart: Strong Scaling GFLOP/sec on strong_summary_all_stencils
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("strong_summary_all_stencils")

$chartObj = $summary.ChartObjects().Add(450, 900, 400, 300)
$chart = $chartObj.Chart
$chart.ChartType = 74

$s1 = $chart.SeriesCollection().NewSeries()
$s1.Name = "n=17"
$s1.XValues = "=n50_strong_alltoallv!`$A`$2:`$A`$12"
$s1.Values = "=n17_strong_overlap_cpu!`$E`$2:`$E`$12"

$s2 = $chart.SeriesCollection().NewSeries()
$s2.Name = "n=31"
$s2.XValues = "=n50_strong_alltoallv!`$A`$2:`$A`$12"
$s2.Values = "=n31_strong_overlap_cpu!`$E`$2:`$E`$12"

$s3 = $chart.SeriesCollection().NewSeries()
$s3.Name = "n=50"
$s3.XValues = "=n50_strong_alltoallv!`$A`$2:`$A`$12"
$s3.Values = "=n50_strong_overlap_cpu!`$E`$2:`$E`$12"

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Strong Scaling GFLOP/sec, N = 4096000 points" + [Environment]::NewLine + "SpMV + MPI (Isend/Irecv Overlap CPU) on Itasca "

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "Process Count (8 ppn)"
$chart.Axes(1).ScaleType = 3
$chart.Axes(1).LogBase = 2

$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "GFLOP/sec"
$chart.Axes(2).ScaleType = 3
$chart.Axes(2).LogBase = 2
$chart.Axes(2).HasMajorGridlines = $true

$chart.HasLegend = $true
$chart.Legend.Position = -4160

# ---------------------------------------------------------------
# 5. Make the summary sheet the active tab (matches activeTab change)
# ---------------------------------------------------------------
$summary.Activate()
